$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the existing "sum" header (G1) onto the new
# "Save" header cell (H1) so it matches the other bold/centered/bordered
# header cells.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Application.CutCopyMode = $false

# Add the new "Save" header and its value for row 2.
$ws.Range("H1").Value = "Save"
$ws.Range("H2").Value = 1
